$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 15.02.2022 09:00"

# Row 8 (Benzina Albert Modřice) gets a fresh price reading.
# The previous reading (Cena/Datum) is shifted into the Old Cena / Old Datum
# columns, and the delta is recorded as text.
$oldCena = $ws.Range("B8").Value2
$newCena = 36.5

$ws.Range("B8").Value2 = $newCena
$ws.Range("C8").Value2 = $oldCena

$delta = [math]::Round($newCena - $oldCena, 2)

# D8/E8 hold text (not numbers/dates) in the refreshed row, so force a text
# format before writing, then clear the format back off so no extra style
# sticks to the cell (matches the "Old Cena"/"Old Datum" columns elsewhere).
$textCells = $ws.Range("D8:E8")
$textCells.NumberFormat = "@"
$ws.Range("D8").Value = [string]$delta
$ws.Range("E8").Value = "2022-02-15 09:03:07"
$textCells.ClearFormats()
